# 10.02/2024 - return to host 12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

$rows = @(
    @("35/65-33",  "ФБел-283",  "42, 30, груз, сер"),
    @("205/55R16", "BEL-262",   "легк, сер, б/к"),
    @("205/55R16", "BEL-317",   "легк, сер, б/к"),
    @("205/55R16", "BEL-317S",  "сер, ошип"),
    @("235/75R15", "BEL-1001",  "легк, сер"),
    @("155/65R13", "BEL-1002",  "легк, сер"),
    @("205/55R16", "BEL-1004",  "легк, сер"),
    @("225/50R17", "BEL-1005",  "легк, сер"),
    @("24.00R35",  "Бел-202",   "210B, Type, сер, C, H"),
    @("24.00R35",  "Бел-212",   "груз, Type, сер, LS-2"),
    @("21.00R35",  "Бел-200",   "202B, Type, сер, C"),
    @("21.00R35",  "Бел-210",   "202B, Type, сер, C, H, LS-2"),
    @("14.00R20",  "BEL-248",   "груз, сер, б/к"),
    @("175/70R13", "Бел-103",   "легк, сер, б/к"),
    @("175/70R13", "Бел-100",   "легк, сер, б/к"),
    @("195/65R15", "Бел-119",   "легк, сер"),
    @("210/80R16", "Бел-777",   "легк, сер"),
    @("215/65R16C","Бел-1000",  "легк, сер"),
    @("205/55R16", "Бел-1001",  "легк, сер"),
    @("225/50R17", "Бел-1005",  "легк, сер")
)

$startRow = 22
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]
    $ws.Cells.Item($r, 5).Value = $item[0]
    $ws.Cells.Item($r, 6).Value = $item[1]
    $ws.Cells.Item($r, 7).Value = $item[2]
    $ws.Cells.Item($r, 8).Value = 2
    $ws.Cells.Item($r, 9).Value = 45335
    $ws.Cells.Item($r, 9).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($r, 10).Value = "нет данных"
}
